# This script updates the "想去人数" (F column) counts on the
# "展览" and "全部类型" worksheets to reflect newly generated data
# (as produced by the gh-pages build at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 48
$ws1.Range("F4").Value  = 617
$ws1.Range("F5").Value  = 156
$ws1.Range("F6").Value  = 9333
$ws1.Range("F7").Value  = 839
$ws1.Range("F9").Value  = 1190
$ws1.Range("F10").Value = 1106
$ws1.Range("F11").Value = 142
$ws1.Range("F12").Value = 82
$ws1.Range("F13").Value = 15
$ws1.Range("F15").Value = 401
$ws1.Range("F16").Value = 85
$ws1.Range("F18").Value = 1240

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 48
$ws4.Range("F6").Value  = 617
$ws4.Range("F7").Value  = 156
$ws4.Range("F8").Value  = 9333
$ws4.Range("F9").Value  = 839
$ws4.Range("F11").Value = 1190
$ws4.Range("F12").Value = 1106
$ws4.Range("F13").Value = 142
$ws4.Range("F14").Value = 82
$ws4.Range("F15").Value = 15
$ws4.Range("F17").Value = 401
$ws4.Range("F18").Value = 85
$ws4.Range("F20").Value = 1240
